$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..G across rows 2..11 (regenerated s_vals data)
$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144)
    3  = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.094976487407548)
    4  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 11.80884686099532)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    6  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    7  = @(0.01514828764759746, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0, 0.7118572775379124)
    8  = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 6.740334628841572)
    9  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 3.645393585217082)
    10 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.997429241610044)
    11 = @(0.04763786555579896, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.369736951971621)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
